# Auto-generated edit script applying the Masamune_Profits.xlsx market-price refresh
# (scheduled runner update) described by the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2394.8572
$ws.Range("I62").Value = 2535.4167
$ws.Range("J62").Value = 1551.5
$ws.Range("K62").Value = 2535.4167
$ws.Range("L62").Value = 1551.5
$ws.Range("M62").Value = -1911.4167
$ws.Range("N62").Value = -2799.5
$ws.Range("H65").Value = 2394.8572
$ws.Range("I65").Value = 2535.4167
$ws.Range("J65").Value = 1551.5
$ws.Range("K65").Value = 12677.0835
$ws.Range("L65").Value = 7757.5
$ws.Range("M65").Value = -9557.083500000001
$ws.Range("N65").Value = -13997.5
$ws.Range("H126").Value = 40399.43
$ws.Range("J126").Value = 40399.43
$ws.Range("L126").Value = 40399.43
$ws.Range("N126").Value = -50279.43
$ws.Range("H128").Value = 53965.2
$ws.Range("J128").Value = 53965.2
$ws.Range("L128").Value = 53965.2
$ws.Range("N128").Value = -63925.2
$ws.Range("H130").Value = 46426.285
$ws.Range("J130").Value = 46426.285
$ws.Range("L130").Value = 46426.285
$ws.Range("N130").Value = -56466.285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1669.68
$ws.Range("I74").Value = 1061
$ws.Range("J74").Value = 2444.3635
$ws.Range("K74").Value = 1061
$ws.Range("L74").Value = 2444.3635
$ws.Range("M74").Value = -187
$ws.Range("N74").Value = -4192.363499999999
$ws.Range("H77").Value = 1669.68
$ws.Range("I77").Value = 1061
$ws.Range("J77").Value = 2444.3635
$ws.Range("K77").Value = 5305
$ws.Range("L77").Value = 12221.8175
$ws.Range("M77").Value = -937
$ws.Range("N77").Value = -20957.8175
$ws.Range("H109").Value = 43377
$ws.Range("J109").Value = 43377
$ws.Range("L109").Value = 43377
$ws.Range("N109").Value = -46151
$ws.Range("H112").Value = 525000000
$ws.Range("J112").Value = 525000000
$ws.Range("L112").Value = 525000000
$ws.Range("N112").Value = -525002954
$ws.Range("H117").Value = 49999.668
$ws.Range("J117").Value = 49999.668
$ws.Range("L117").Value = 49999.668
$ws.Range("N117").Value = -59177.668
$ws.Range("H123").Value = 51000
$ws.Range("J123").Value = 51000
$ws.Range("L123").Value = 51000
$ws.Range("N123").Value = -60800

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -59178

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6366.0464
$ws.Range("I31").Value = 3467.9
$ws.Range("J31").Value = 7244.273
$ws.Range("K31").Value = 3467.9
$ws.Range("L31").Value = 7244.273
$ws.Range("M31").Value = -3172.9
$ws.Range("N31").Value = -7834.273
$ws.Range("H34").Value = 6366.0464
$ws.Range("I34").Value = 3467.9
$ws.Range("J34").Value = 7244.273
$ws.Range("K34").Value = 3467.9
$ws.Range("L34").Value = 7244.273
$ws.Range("M34").Value = -3265.9
$ws.Range("N34").Value = -7648.273
$ws.Range("H100").Value = 46996
$ws.Range("J100").Value = 46996
$ws.Range("L100").Value = 46996
$ws.Range("N100").Value = -49160
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H118").Value = 44734
$ws.Range("J118").Value = 44734
$ws.Range("L118").Value = 44734
$ws.Range("N118").Value = -48048

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 200003180
$ws.Range("I80").Value = 9999
$ws.Range("J80").Value = 250001470
$ws.Range("K80").Value = 29997
$ws.Range("L80").Value = 750004410
$ws.Range("M80").Value = -29061
$ws.Range("N80").Value = -750006282
$ws.Range("H83").Value = 200003180
$ws.Range("I83").Value = 9999
$ws.Range("J83").Value = 250001470
$ws.Range("K83").Value = 89991
$ws.Range("L83").Value = 2250013230
$ws.Range("M83").Value = -85311
$ws.Range("N83").Value = -2250022590

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2070
$ws.Range("I122").Value = 2316
$ws.Range("J122").Value = 1933.3334
$ws.Range("K122").Value = 6948
$ws.Range("L122").Value = 5800.0002
$ws.Range("M122").Value = -4498
$ws.Range("N122").Value = -10700.0002
$ws.Range("H133").Value = 27472.957
$ws.Range("J133").Value = 27472.957
$ws.Range("L133").Value = 27472.957
$ws.Range("N133").Value = -37592.95699999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 15000
$ws.Range("J54").Value = 15000
$ws.Range("L54").Value = 15000
$ws.Range("N54").Value = -16288
$ws.Range("H61").Value = 2600.2144
$ws.Range("I61").Value = 2600.2144
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2600.2144
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2398.2144
$ws.Range("N61").ClearContents()
$ws.Range("H81").Value = 42181
$ws.Range("J81").Value = 42181
$ws.Range("L81").Value = 42181
$ws.Range("N81").Value = -44177
$ws.Range("H84").Value = 42181
$ws.Range("J84").Value = 42181
$ws.Range("L84").Value = 126543
$ws.Range("N84").Value = -136527
$ws.Range("H110").Value = 40644
$ws.Range("J110").Value = 40644
$ws.Range("L110").Value = 40644
$ws.Range("N110").Value = -48824
$ws.Range("H113").Value = 2600.2144
$ws.Range("I113").Value = 2600.2144
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2600.2144
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -430.2143999999998
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 2742.8572
$ws.Range("I122").Value = 2666.6667
$ws.Range("J122").Value = 2800
$ws.Range("K122").Value = 8000.000100000001
$ws.Range("L122").Value = 8400
$ws.Range("M122").Value = -5550.000100000001
$ws.Range("N122").Value = -13300
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 664.2
$ws.Range("I113").Value = 597.1539
$ws.Range("J113").Value = 1100
$ws.Range("K113").Value = 1791.4617
$ws.Range("L113").Value = 3300
$ws.Range("M113").Value = 378.5382999999999
$ws.Range("N113").Value = -7640
$ws.Range("H115").Value = 37371.668
$ws.Range("J115").Value = 37371.668
$ws.Range("L115").Value = 37371.668
$ws.Range("N115").Value = -40505.668
$ws.Range("H116").Value = 47654
$ws.Range("J116").Value = 47654
$ws.Range("L116").Value = 47654
$ws.Range("N116").Value = -56832
$ws.Range("H118").Value = 31517
$ws.Range("J118").Value = 31517
$ws.Range("L118").Value = 31517
$ws.Range("N118").Value = -34831
$ws.Range("H119").Value = 47273.43
$ws.Range("J119").Value = 47273.43
$ws.Range("L119").Value = 47273.43
$ws.Range("N119").Value = -56949.43
$ws.Range("H120").Value = 42424.332
$ws.Range("J120").Value = 42424.332
$ws.Range("L120").Value = 42424.332
$ws.Range("N120").Value = -52100.332
